# Swap the order of recorded-by names in column G:
# "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#
# We use Find/FindNext (rather than scanning every cell) so that empty
# cells are never touched/materialized.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$addresses = New-Object System.Collections.ArrayList

$first = $ws.Cells.Find($oldValue)
if ($first -ne $null) {
    $firstAddr = $first.Address()
    $current = $first
    $count = 0
    do {
        [void]$addresses.Add($current.Address())
        $count = $count + 1
        $current = $ws.Cells.FindNext($current)
    } while ($current -ne $null -and $current.Address() -ne $firstAddr -and $count -lt 1000)
}

foreach ($addr in $addresses) {
    $ws.Range($addr).Value = $newValue
}
